$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# --- Paragraph formatting changes on paragraph 1 -----------------------
# Give paragraph 1 the same paragraph border (5pt space on all sides, no
# visible line) that the other body paragraphs already have.
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Left indent: 120 twips (6pt) -> 225 twips (11.25pt)
$p1.Format.LeftIndent = 11.25

# --- Text/run changes on paragraph 1 ------------------------------------
# Paragraph 1 originally holds two runs:
#   run A: "**ID__AFFARS_pgi_5335_topic_8__ID**"
#   run B: " "   (a single trailing space)
# We need run A's text changed and run B removed entirely (not just
# emptied), so locate each run's range via Find and edit precisely.

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("**ID__AFFARS_pgi_5335_topic_8__ID**")
if ($found) {
    $idRange = $d.Range($findRange.Start, $findRange.End)

    # The trailing single-space run immediately follows the ID run.
    $spaceRange = $d.Range($idRange.End, $idRange.End + 1)
    if ($spaceRange.Text -eq " ") {
        $spaceRange.Text = ""
    }

    # Now replace the ID run's own text (range length unaffected by the
    # deletion above, since it came after this range).
    $idRange.Text = "**ID__AFFARS_SMC_PGI_5335_017_3__ID**"
}
